$d = $word.ActiveDocument

# The document currently has a single inline picture (a 1x1 placeholder
# image, embedded via relationship rId23) sitting right after the
# "Terrace-I" bookmark. The edit replaces that picture with a visible
# hyperlink run (styled with the built-in "Hyperlink" character style)
# that points to the image's real URL on ura.gov.sg, reusing the same
# relationship id (rId23) the picture used to reference.

$url = "https://ura.gov.sg/-/media/Corporate/Guidelines/Development-control/Landed-Housing/TH01_Plot_Size_Width_Terrace_1.jpg?h=100%25&w=100%25"

# Find the inline picture anywhere in the document (there is exactly one).
$shape = $null
for ($i = 1; $i -le $d.InlineShapes.Count; $i++) {
    $candidate = $d.InlineShapes.Item($i)
    if ($candidate.Type -eq 3) {
        $shape = $candidate
        break
    }
}
if ($shape -eq $null -and $d.InlineShapes.Count -gt 0) {
    $shape = $d.InlineShapes.Item(1)
}

$target = $shape.Range
$shape.Delete()

$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships">
        <w:body>
          <w:p>
            <w:hyperlink r:id="rId23">
              <w:r>
                <w:rPr>
                  <w:rStyle w:val="Hyperlink"/>
                </w:rPr>
                <w:t xml:space="preserve">https://ura.gov.sg/-/media/Corporate/Guidelines/Development-control/Landed-Housing/TH01_Plot_Size_Width_Terrace_1.jpg?h=100%25&amp;w=100%25</w:t>
              </w:r>
            </w:hyperlink>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$target.InsertXML($xml)

Write-Host "Replaced inline image with hyperlink run pointing to: $url"
